$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (12) to the talk/venue grid:
#   A12 = talk title, D12 = date for the "ThatConference" column
$ws.Range("A12").Value = "Async Masterclass"
$ws.Range("D12").Value = "2023-07"

# Move the active selection to A13, matching the post-edit state
$ws.Range("A13").Select()
